$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill the whole of row 4 (columns C through AD) with "Y" to match the
# extended range of the template (previously only a sparse subset of
# columns C:Z had "Y" placeholders).
$ws.Range("C4:AD4").Value = "Y"

# New header in AD3: "Other"
$ws.Range("AD3").Value = "Other"

# Update the active selection to L8 (as captured in the saved view state).
$ws.Range("L8").Select()
